# Apply odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3

# Row 6
$ws.Range("H6").Value = 3.45
$ws.Range("I6").Value = 5.4
$ws.Range("J6").Value = 2.18
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 6.25
$ws.Range("S6").Value = 1.42
$ws.Range("T6").Value = 2.45
$ws.Range("Z6").Value = 11.75
$ws.Range("AB6").Value = 35
$ws.Range("AC6").Value = 7.8
$ws.Range("AH6").Value = 12
$ws.Range("AI6").Value = 32
$ws.Range("AO6").Value = 7.8
$ws.Range("AP6").Value = 19.5
$ws.Range("AR6").Value = 65
$ws.Range("AS6").Value = 300
$ws.Range("AT6").Value = 2.42
$ws.Range("AU6").Value = 8
$ws.Range("AW6").Value = 6.8
$ws.Range("BA6").Value = 250

# Row 7
$ws.Range("N7").Value = 8
$ws.Range("S7").Value = 1.35
$ws.Range("T7").Value = 2.94
